$wb = $excel.ActiveWorkbook

# Overview sheet: the 9a3c42df... file row (row 3) moved from "Ready for
# handoff" to "Handed back: in sync with en-US" for both the zh-cn and
# de-de status columns now that the handback report has been generated.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn detail sheet: same row's Status flips to "Handed back", and the
# Latest Handback DateTime is stamped with the handback report time.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-20 00:37:16"

# de-de detail sheet: same change, with its own handback timestamp.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-20 00:37:21"
